$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.852.44'
$ws.Cells.Item(2, 5).Value = '  +0.37%  '
$ws.Cells.Item(3, 4).Value = '2.658.15'
$ws.Cells.Item(3, 5).Value = '  +4.31%  '
$ws.Cells.Item(4, 5).Value = '  +0.13%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '514.10'
$ws.Cells.Item(5, 4).NumberFormat = 'General'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.99%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '144.08'
$ws.Cells.Item(6, 4).NumberFormat = 'General'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.74%  '
$ws.Cells.Item(7, 5).Value = '  -0.50%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.567'
$ws.Cells.Item(8, 4).NumberFormat = 'General'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +3.17%  '
$ws.Cells.Item(9, 4).Value = '2.692.20'
$ws.Cells.Item(9, 5).Value = '  +5.62%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '6.23'
$ws.Cells.Item(10, 4).NumberFormat = 'General'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +0.88%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.106'
$ws.Cells.Item(11, 4).NumberFormat = 'General'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +5.95%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.335'
$ws.Cells.Item(12, 4).NumberFormat = 'General'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +2.29%  '
$ws.Cells.Item(13, 5).Value = '  -0.92%  '
$ws.Cells.Item(14, 4).Value = '3.134.24'
$ws.Cells.Item(14, 5).Value = '  +4.42%  '
$ws.Cells.Item(15, 4).Value = '58.874.88'
$ws.Cells.Item(15, 5).Value = '  +0.40%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '21.00'
$ws.Cells.Item(16, 4).NumberFormat = 'General'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +2.63%  '
$ws.Cells.Item(17, 5).Value = '  +2.79%  '
$ws.Cells.Item(18, 4).Value = '2.686.87'
$ws.Cells.Item(18, 5).Value = '  +5.22%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '347.57'
$ws.Cells.Item(19, 4).NumberFormat = 'General'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +5.55%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '4.53'
$ws.Cells.Item(20, 4).NumberFormat = 'General'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +1.11%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '10.41'
$ws.Cells.Item(21, 4).NumberFormat = 'General'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +4.27%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.14'
$ws.Cells.Item(22, 4).NumberFormat = 'General'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +4.16%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.998'
$ws.Cells.Item(23, 4).NumberFormat = 'General'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.21%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '60.95'
$ws.Cells.Item(24, 4).NumberFormat = 'General'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +3.03%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.419'
$ws.Cells.Item(25, 4).NumberFormat = 'General'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +3.81%  '
$ws.Cells.Item(26, 4).Value = '2.787.36'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.993'
$ws.Cells.Item(27, 4).NumberFormat = 'General'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.65%  '
$ws.Cells.Item(28, 5).Value = '  +2.43%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0809'
$ws.Cells.Item(29, 5).Value = '  +5.31%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.26'
$ws.Cells.Item(30, 4).NumberFormat = 'General'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +6.46%  '
$ws.Cells.Item(31, 5).Value = '  -0.40%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '6.47'
$ws.Cells.Item(32, 4).NumberFormat = 'General'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +11.97%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '18.98'
$ws.Cells.Item(33, 4).NumberFormat = 'General'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +2.88%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.57'
$ws.Cells.Item(34, 4).NumberFormat = 'General'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +2.86%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '150.00'
$ws.Cells.Item(35, 4).NumberFormat = 'General'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.54%  '
$ws.Cells.Item(36, 5).Value = '  +15.15%  '
$ws.Cells.Item(37, 5).Value = '  +3.41%  '
$ws.Cells.Item(38, 5).Value = '  +4.12%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '36.77'
$ws.Cells.Item(39, 4).NumberFormat = 'General'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +3.14%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.844'
$ws.Cells.Item(40, 4).NumberFormat = 'General'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +3.04%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.68'
$ws.Cells.Item(41, 4).NumberFormat = 'General'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +6.33%  '
$ws.Cells.Item(42, 5).Value = '  +2.06%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.621'
$ws.Cells.Item(43, 4).NumberFormat = 'General'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +3.05%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '278.34'
$ws.Cells.Item(44, 4).NumberFormat = 'General'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.51%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '19.82'
$ws.Cells.Item(46, 4).NumberFormat = 'General'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +7.47%  '
$ws.Cells.Item(47, 2).Value = 'Stellar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0981'
$ws.Cells.Item(47, 4).NumberFormat = 'General'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.68%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0533'
$ws.Cells.Item(48, 4).NumberFormat = 'General'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.30%  '
$ws.Cells.Item(49, 5).Value = '  +2.50%  '
$ws.Cells.Item(50, 4).Value = '2.006.59'
$ws.Cells.Item(50, 5).Value = '  +5.65%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '4.73'
$ws.Cells.Item(51, 4).NumberFormat = 'General'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +5.06%  '
